$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B-E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (columns B-E)
$ws.Range("B2").Value = 11.27668584925304
$ws.Range("C2").Value = 10.937389317996539
$ws.Range("D2").Value = 12.296540643036979
$ws.Range("E2").Value = 11.148224681782706

# Row 3 data values (columns B-E)
$ws.Range("B3").Value = 10.67466473685695
$ws.Range("C3").Value = 9.050227224310861
$ws.Range("D3").Value = 10.630007628015582
$ws.Range("E3").Value = 10.856173603842553

# Update the selection to match the edited range
$ws.Range("B1:E3").Select()
